# Regenerate merged AHB files
# - rename header columns: "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410"
# - wrap the data range in a table (Table1)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# 1) Rename the "_old" / "_new" header cells to "_FV2404" / "_FV2410"
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2404")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2410")
}

# 2) Turn the used range into an Excel Table ("Table1") with a header row
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the top (header) row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
